$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 348; this shifts the existing rows 348..390
# down to 349..391 (matching the target diff where all old rows 348-390
# reappear one row lower, and a brand new data row lands at 348).
$ws.Rows.Item(348).Insert()

# Populate the newly inserted row 348 with the new weekly price entry.
$ws.Cells.Item(348, 1).Value = 7
$ws.Cells.Item(348, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(348, 3).Value = "Ñuble"
$ws.Cells.Item(348, 4).Value = 45212
$ws.Cells.Item(348, 5).Value = 16
$ws.Cells.Item(348, 6).Value = 100112045
$ws.Cells.Item(348, 7).Value = "Zapallo"
$ws.Cells.Item(348, 8).Value = "Paine"
$ws.Cells.Item(348, 9).Value = "1a (guarda)"
$ws.Cells.Item(348, 10).Value = 200
$ws.Cells.Item(348, 11).Value = 450
$ws.Cells.Item(348, 12).Value = 450
$ws.Cells.Item(348, 13).Value = 450
$ws.Cells.Item(348, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(348, 15).Value = "Región del Maule"
$ws.Cells.Item(348, 16).Value = 450
$ws.Cells.Item(348, 17).Value = 1
$ws.Cells.Item(348, 18).Value = "Hortaliza"
